$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44281
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 5500
$ws.Range("P2").Value = 5500

# Row 3
$ws.Range("D3").Value = 44259
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = 4250
$ws.Range("P3").Value = 4250

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("P4").Value = 8500

# Row 5
$ws.Range("D5").Value = 44253
$ws.Range("H5").Value = "Americana (o)"
$ws.Range("I5").Value = "Segunda"
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = 4250
$ws.Range("P5").Value = 4250

# Row 6
$ws.Range("D6").Value = 44371
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 7375
$ws.Range("P6").Value = 7375

# Row 7
$ws.Range("D7").Value = 44410
$ws.Range("K7").Value = 5500
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5750
$ws.Range("P7").Value = 5750

# Row 8
$ws.Range("D8").Value = 44263
$ws.Range("J8").Value = 100
$ws.Range("M8").Value = 7500
$ws.Range("P8").Value = 7500

# Row 9
$ws.Range("D9").Value = 44414
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6500
$ws.Range("P9").Value = 6500
